$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$BValues = @(1.018308750443737, 0.9259691007000583, 0.8696061280359686, 0.8467219209670418, 0.8429271048618716, 0.8692971625923178, 0.9864008094093606, 1.218694951279929, 1.391004657586222, 1.469756287192808, 1.499630411409441, 1.493194148223665, 1.47221299644383, 1.459368278486352, 1.385865423127996, 1.340867678678819, 1.315020728531408, 1.306275337519708, 1.345654189762456, 1.478374243502742, 1.565420929560446, 1.518934495689166, 1.343490136906439, 1.155566517308557)
$CValues = @(0.1826043342624644, 0.1603048161715037, 0.1466179368479175, 0.141041733960094, 0.1401158897575101, 0.1465427288773355, 0.1749144065923645, 0.2305935653768358, 0.2715336225799376, 0.2901671267158008, 0.2972245905245927, 0.2957045793743305, 0.2907477214124867, 0.2877116790768071, 0.2703160711739656, 0.2596468830479637, 0.2535111597438231, 0.2514338667643869, 0.260782542218891, 0.2922036340634691, 0.3127470500788831, 0.3017819272378972, 0.2602691165168665, 0.215525652017476)
$DValues = @(0.1936504619927248, 0.1926493416014239, 0.1920882087606017, 0.1918730617824238, 0.1918381548338317, 0.1920852524057395, 0.1932941895668137, 0.196087886439841, 0.1983959808048468, 0.1995010315099677, 0.1999273621431144, 0.1998351948918753, 0.1995359485053569, 0.1993536751458436, 0.1983248673894664, 0.1977078008730402, 0.1973580685173459, 0.197240547707942, 0.1977729520801361, 0.1996236311743615, 0.2008790178147279, 0.2002048143143469, 0.1977434815564791, 0.1952870735939953)
$FValues = @(1.691700264227777, 1.679769438616901, 1.673525869143148, 1.671253056808467, 1.670892038528834, 1.673494118718281, 1.687361582691921, 1.723171300981974, 1.754780937684302, 1.770322641967496, 1.776375793934378, 1.775064663061485, 1.770817270479299, 1.768237499409125, 1.753788699800097, 1.745223027639142, 1.740405640283825, 1.738793317107678, 1.74612353458113, 1.772060271697043, 1.789990094784926, 1.780330828034252, 1.745716081742074, 1.712556529289856)
$GValues = @(1.034880089528542, 1.0206957595001, 1.01277849485524, 1.009750524098834, 1.009259691051597, 1.012736856374701, 1.029824487974082, 1.069655239959303, 1.102830385954633, 1.11878469278335, 1.124951180876764, 1.123617546699506, 1.119289504497118, 1.116654750354058, 1.101805170109657, 1.092917028756176, 1.087885960084662, 1.086196438611637, 1.0938547810307, 1.120557357626581, 1.138737680491829, 1.12896753893807, 1.093430577709768, 1.058197236135015)
$HValues = @(1.007899836459046, 1.007137010388632, 1.007250615564217, 1.00744303071059, 1.007483799890366, 1.007252619210291, 1.007515884605084, 1.012661557763835, 1.019282588966576, 1.0229155890801, 1.024380922365538, 1.024061347249244, 1.023034345656328, 1.022416954024322, 1.019057685662119, 1.017156155010895, 1.016120884592453, 1.015780388187096, 1.017352525854903, 1.023333567029397, 1.027764887146645, 1.025351916045935, 1.017263566123034, 1.010772028089008)
$JValues = @(0.2507695397009115, 0.2506886421895089, 0.2507595185482998, 0.2508186951601559, 0.2508303505342226, 0.2507601939898265, 0.2507166144476471, 0.2515888197819578, 0.2528155404039012, 0.2535013250250131, 0.2537794162298539, 0.2537187055057615, 0.2535238348284139, 0.2534068679692183, 0.2527732958470068, 0.2524173592060635, 0.2522246558248398, 0.2521614736465381, 0.2524540048270865, 0.2535805734708489, 0.2544241117597892, 0.2539640750200149, 0.2524374001864373, 0.2512501265865765)
$MValues = @(0.4514141866378054, 0.4271570723711591, 0.4124728992235731, 0.4065417942726057, 0.4055601305827352, 0.4123926963398574, 0.4430068181392386, 0.5047066609107134, 0.5510598152926605, 0.5723708980975033, 0.5804732157810975, 0.5787268021504985, 0.5730368341394296, 0.569555765838075, 0.5496716048777728, 0.5375308543181916, 0.5305689779500184, 0.5282154425316463, 0.538821068122644, 0.5747072401931774, 0.598349018422013, 0.5857137628130147, 0.5382377063826027, 0.4878362026962222)
$NValues = @(1.383031793137036, 1.401937714846075, 1.414124332706895, 1.419235734976677, 1.420093247872558, 1.414192678948536, 1.389430430854066, 1.345466908601774, 1.315976231907537, 1.303171941831053, 1.298411371017744, 1.299432721191321, 1.302778518972279, 1.304839402167141, 1.31682534193135, 1.324335038900271, 1.32871192852755, 1.330203738082798, 1.323529664936074, 1.301793383047289, 1.288101310577376, 1.295361926172806, 1.323893589293851, 1.356867458240956)

for ($i = 0; $i -lt 24; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, "B").Value = $BValues[$i]
    $ws.Cells.Item($row, "C").Value = $CValues[$i]
    $ws.Cells.Item($row, "D").Value = $DValues[$i]
    $ws.Cells.Item($row, "F").Value = $FValues[$i]
    $ws.Cells.Item($row, "G").Value = $GValues[$i]
    $ws.Cells.Item($row, "H").Value = $HValues[$i]
    $ws.Cells.Item($row, "J").Value = $JValues[$i]
    $ws.Cells.Item($row, "M").Value = $MValues[$i]
    $ws.Cells.Item($row, "N").Value = $NValues[$i]
}

Write-Host "applied 380 kV case values"